# Update the "想去人数" (interested-count) figures on both the
# "展览" and "全部类型" sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1051
$ws1.Range("F4").Value = 175
$ws1.Range("F5").Value = 2881
$ws1.Range("F7").Value = 266
$ws1.Range("F11").Value = 130
$ws1.Range("F12").Value = 50
$ws1.Range("F13").Value = 2709
$ws1.Range("F14").Value = 961

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1051
$ws4.Range("F5").Value = 175
$ws4.Range("F6").Value = 2881
$ws4.Range("F8").Value = 266
$ws4.Range("F13").Value = 130
$ws4.Range("F14").Value = 50
$ws4.Range("F15").Value = 2709
$ws4.Range("F16").Value = 961
